$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: additional sub-task for the exploratory-analysis group
$ws.Range("D16").Value = "*tendencias de tasa de los municipios con tasa más alta"

# New row 24: new task under "transformaciones de la data:"
$ws.Range("D24").Value = "agregar las variables a los features ya creados"
$ws.Range("E24").Value = "felipe"

# Rows 33-34: assign "esteban" to existing cluster tasks
$ws.Range("E33").Value = "esteban"
$ws.Range("E34").Value = "esteban"

# New row 36: new task
$ws.Range("D36").Value = "modelo predicción casos mensuales"
$ws.Range("E36").Value = "esteban"

# Row 11: update assignee text (was "Joha, David, Daniel")
$ws.Range("E11").Value = "David, Daniel"

# Rows 12-16: add same assignee to sub-tasks under "exploración de la datam /exploratory Analysis"
$ws.Range("E12").Value = "David, Daniel"
$ws.Range("E13").Value = "David, Daniel"
$ws.Range("E14").Value = "David, Daniel"
$ws.Range("E15").Value = "David, Daniel"
$ws.Range("E16").Value = "David, Daniel"

# Update selection/active cell
$ws.Range("D11").Select()
